$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Add two new rows to the attrition table (Table3)
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# Row 7: Department Leader
$ws.Range("A2").Copy($ws.Range("A7"))
$ws.Range("B7").Value = "Department Leader"
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 2
$ws.Range("F6").Copy($ws.Range("F7"))

# Row 8: Regional Leader
$ws.Range("A2").Copy($ws.Range("A8"))
$ws.Range("B8").Value = "Regional Leader"
$ws.Range("C8").Value = 11
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 11
$ws.Range("F6").Copy($ws.Range("F8"))
